# Delete the record with Kayıt No 11362865 from both the "Kayitlar"
# summary sheet and the "Merkez İlçe" district sheet. Deleting the row
# shifts every subsequent row up by one, shrinking the used range by a
# single row on each sheet (matches the commit message:
# "Kayıt silindi: 11362865").

$wb = $excel.ActiveWorkbook
$recordId = "11362865"

foreach ($sheetName in @("Kayitlar", "Merkez İlçe")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $match = $ws.Columns.Item(1).Find($recordId)
    if ($match -ne $null) {
        $ws.Rows.Item($match.Row).Delete()
    }
}
